$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/tab
$ws.Name = "GammaFiber2F"

# Tiny precision corrections in row 13 (results of the new Gaussian
# Quadrature based averaging scheme)
$ws.Range("D13").Value = 0.9980709158298056
$ws.Range("I13").Value = 0.9987901497843227
$ws.Range("J13").Value = 0.9980709158298056
$ws.Range("K13").Value = 0.9992227491542817
$ws.Range("L13").Value = 0.9974658446126285

# Append a new row (14 -> HexGrid-60degTilt5degRes) with the averaged
# intensities computed by the new Gaussian Quadrature scheme
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16").Value = 1.186928358257481
$ws.Range("D16").Value = 0.6025458675910295
$ws.Range("E16").Value = 1.05201901139688
$ws.Range("F16").Value = 1.186928358257481
$ws.Range("G16").Value = 0.7990077946846026
$ws.Range("H16").Value = 1.1312998136103
$ws.Range("I16").Value = 1.09943573977727
$ws.Range("J16").Value = 0.6025458675910295
$ws.Range("K16").Value = 0.8272824394939546
$ws.Range("L16").Value = 1.007105398875717
$ws.Range("M16").Value = 0.9785394308862604

# Match the formatting (bold, centered, bordered) already used for the
# other cells in column A by copying it down from the row above
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$excel.CutCopyMode = $false
